$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# This script regenerates the localization-status report: the row that
# used to describe "dcdaf9f5-...md" has now been handed off ("Ready for
# handoff") while "ab3f5bcd-...md" finished its handback cycle and is
# "Handed back: in sync with en-US". The generator re-emits the sheets
# with dcdaf9f5 first (row 2) and ab3f5bcd second (row 3).
# ----------------------------------------------------------------------

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-35-19 00:35:59"

$ov.Range("A3").Value = "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-36-19 00:36:50"

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf.c51fdbc4c7c41c142e31c0a58ae03efdf802fcf2.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-19 00:35:56"
$zh.Range("F2").Value = "dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf.md"
$zh.Range("G2").Value = "dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf.c51fdbc4c7c41c142e31c0a58ae03efdf802fcf2.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-19 00:36:19"
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-19 00:36:48"
$zh.Range("F3").Value = "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md"
$zh.Range("G3").Value = "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-19 00:36:19"
$zh.Range("I3").Value = "Include"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf.c51fdbc4c7c41c142e31c0a58ae03efdf802fcf2.de-de.xlf"
$de.Range("E2").Value = "2016-03-19 00:35:59"
$de.Range("F2").Value = "dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf.md"
$de.Range("G2").Value = "dcdaf9f5-af2c-4b5f-83bf-f8f5a6c29cdf.c51fdbc4c7c41c142e31c0a58ae03efdf802fcf2.de-de.xlf"
$de.Range("H2").Value = "2016-03-19 00:36:25"
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.de-de.xlf"
$de.Range("E3").Value = "2016-03-19 00:36:50"
$de.Range("F3").Value = "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.md"
$de.Range("G3").Value = "ab3f5bcd-c748-4009-8642-b2e8ef5797d8.f3d8254288d08074997ef0217cd21a9d62cc2cab.de-de.xlf"
$de.Range("H3").Value = "2016-03-19 00:36:25"
$de.Range("I3").Value = "Include"
